$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weather station ("Honolulu:") was added to the list, inserted as a
# new row between "Anchorage:" (row 17) and "Corpus Christi:" (old row 18).
$ws.Rows.Item(18).Insert() | Out-Null

$ws.Range("A18").Value = "Honolulu:"
$ws.Range("B18").Value = "USW00022521"

# The inserted row picks up row 17's formatting for the trailing (empty)
# columns; put C18:D18 back to the plain/default look used by the other
# blank filler cells in this column.
$ws.Range("C18:D18").ClearFormats() | Out-Null

# Reflect the final cursor/selection position left behind by the edit.
$ws.Range("D17").Select() | Out-Null
